$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------------
# 1) Insert a new (still empty) row at 22, pushing the existing
#    "Toxic Spores / Blast" row down to row 23. We fill in its values last
#    so that the new shared strings end up in the same relative order as
#    the authoring workbook.
# --------------------------------------------------------------------------
$ws.Rows.Item(22).Insert()

# Bring over the formatting of the (now shifted) original row so the new
# row 22 matches the table's look.
$ws.Range("B23:L23").Copy()
$ws.Range("B22:L22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A22").Clear()

# --------------------------------------------------------------------------
# 2) Two brand-new skill rows get appended after the existing data (rows 24
#    and 25): "Bettle / Assault spikes" and "Fire Bettle / Magma jet".
#    Column C is entered before column B for each row (matches the shared
#    string ordering of the saved workbook).
# --------------------------------------------------------------------------

# Row 24 - Bettle / Assault spikes
$ws.Range("C24").Value = "Assault spikes"
$ws.Range("B24").Value = "Bettle"
$ws.Range("D24").Value = "switch"
$ws.Range("E24").Value = -1
$ws.Range("F24").Value = "self"
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 1
$ws.Range("I24").Value = "armor/damage/range"
$ws.Range("J24").Value = "2/4/-2"
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = "0"

# Row 25 - Fire Bettle / Magma jet
$ws.Range("C25").Value = "Magma jet"
$ws.Range("B25").Value = "Fire Bettle"
$ws.Range("D25").Value = "damage"
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = "a"
$ws.Range("G25").Value = 50
$ws.Range("H25").Value = 1
$ws.Range("I25").Value = "n"
$ws.Range("J25").Value = "75"
$ws.Range("K25").Value = 6
$ws.Range("L25").Value = "1"

# Column C is the (bold, bordered) "skill name" style, normally used on
# column B; swap it onto C24/C25 to match.
$ws.Range("B15").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("C25").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C24").Value = "Assault spikes"
$ws.Range("C25").Value = "Magma jet"

# --------------------------------------------------------------------------
# 3) Fill in the new row 22 - "Toxic Spores / Charge"
# --------------------------------------------------------------------------
$ws.Range("C22").Value = "Charge"
$ws.Range("B22").Value = "Toxic Spores"
$ws.Range("D22").Value = "continuous"
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = "self"
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 4
$ws.Range("I22").Value = "speed"
$ws.Range("J22").Value = "3"
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = "0"

# --------------------------------------------------------------------------
# 4) Update the view: scrolled down a bit, selection on M25.
# --------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("M25").Select()
